$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(329).Insert()

$ws.Range("A329").Value = 3
$ws.Range("B329").Value = "Femacal de La Calera"
$ws.Range("C329").Value = "Coquimbo"
$ws.Range("D329").Value = 44876
$ws.Range("E329").Value = 5
$ws.Range("F329").Value = 100112009
$ws.Range("G329").Value = "Acelga"
$ws.Range("H329").Value = "Sin especificar"
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 230
$ws.Range("K329").Value = 3500
$ws.Range("L329").Value = 4000
$ws.Range("M329").Value = 3739
$ws.Range("N329").Value = '$/docena de atados (6 kilos)'
$ws.Range("O329").Value = "Provincia de Quillota"
$ws.Range("P329").Value = 623
$ws.Range("Q329").Value = 6
$ws.Range("R329").Value = "Hortaliza"
